# Updated ping statistics in the "sokhnaReport" sheet (alex report script).
# The monitoring run now produced a handful of dropped/timed-out pings for
# several hosts, so the Online / Not Reachable / Timeout percentage columns
# (C, D, E) need to be refreshed to reflect the new numbers, and column E's
# width needs to shrink now that the values are shorter.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new percentage value, taken from the updated report.
$updates = @{
    "C2" = 99.83;
    "D2" = 0.08;
    "E2" = 0.08;
    "C3" = 99.83;
    "D3" = 0.08;
    "E3" = 0.08;
    "C4" = 99.81;
    "D4" = 0.1;
    "E4" = 0.1;
    "C5" = 99.88;
    "D5" = 0.06;
    "E5" = 0.06;
    "C6" = 99.88;
    "D6" = 0.06;
    "E6" = 0.06;
    "D8" = 0.4;
    "E8" = 99.6;
    "C9" = 99.86;
    "D9" = 0.07;
    "E9" = 0.07;
    "D10" = 0.45;
    "E10" = 99.55;
    "C13" = 99.87;
    "D13" = 0.06;
    "E13" = 0.06;
    "C15" = 99.87;
    "D15" = 0.06;
    "E15" = 0.06;
    "C17" = 99.85;
    "D17" = 0.08;
    "E17" = 0.08;
    "C19" = 99.81;
    "D19" = 0.1;
    "E19" = 0.1;
    "C21" = 99.83;
    "D21" = 0.08;
    "E21" = 0.08;
    "C22" = 99.83;
    "D22" = 0.08;
    "E22" = 0.08;
    "C29" = 89.37;
    "D29" = 10.55;
    "E29" = 0.08;
    "C30" = 89.4;
    "D30" = 10.54;
    "E30" = 0.06;
    "C31" = 89.37;
    "D31" = 10.55;
    "E31" = 0.08;
    "C33" = 89.39;
    "D33" = 10.54;
    "E33" = 0.07;
    "C35" = 89.41;
    "D35" = 10.53;
    "E35" = 0.06;
    "C36" = 89.41;
    "D36" = 10.53;
    "E36" = 0.06;
    "C37" = 89.42;
    "D37" = 10.52;
    "E37" = 0.05;
    "C38" = 89.4;
    "D38" = 10.54;
    "E38" = 0.06;
    "C39" = 89.36;
    "D39" = 10.54;
    "E39" = 0.1;
    "C40" = 25.54;
    "D40" = 24.9;
    "E40" = 49.56;
    "C41" = 89.14;
    "D41" = 10.6;
    "E41" = 0.26;
    "C42" = 99.85;
    "D42" = 0.08;
    "E42" = 0.08;
    "C43" = 21.02;
    "D43" = 78.96;
    "E43" = 0.02;
    "C44" = 13.32;
    "D44" = 22.27;
    "E44" = 64.41;
    "C46" = 12.93;
    "D46" = 22.4;
    "E46" = 64.68;
    "C47" = 9.59;
    "D47" = 22.2;
    "E47" = 68.22
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# Column E (Timeout %) is narrower now that the long repeating-decimal
# values (e.g. 66.6666666666667) have been replaced with short ones.
$ws.Columns.Item(5).ColumnWidth = 7.666666666666667
